# ---------------------------------------------------------------------------
# Reproduces the commit:
#   "Removed some magic numbers; added downloading newest data feature;
#    started prep for visualization more than temperature; made pre-PoC
#    of web app"
#
# Concretely, against the "bands" workbook this means:
#   1. On sheet "bands": widen column B a bit, and colour-code a handful of
#      row-number cells in column A - orange (FFC000) for rows that were
#      "pulled out" into a new summary sheet, plain white for a few rows
#      that had a highlight cleared.
#   2. Add a new sheet "Arkusz1" (after "bands") that collects the rows that
#      got the orange highlight into a small standalone table, and make it
#      the active / selected sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$bands = $wb.Worksheets.Item("bands")

# --- 1a. Column B is a bit wider now ---------------------------------------
$bands.Columns.Item(2).ColumnWidth = 37.3

# --- 1b/1c. Highlight the "moved out" rows in column A with orange --------
# (first cell touched allocates the new orange fill / cellXf)
$orangeRows = 11,146,147,415,417,420,421,424,432
foreach ($r in $orangeRows) {
    $bands.Cells.Item($r, 1).Interior.Color = 49407
}

# Row 433 gets the same orange fill, but picked up a (cosmetically
# identical) "applied font" flag along the way in the original edit.
$c433 = $bands.Cells.Item(433, 1)
$c433.Interior.Color = 49407
$c433.Font.ThemeColor = 1

# Remaining plain-orange rows
$bands.Cells.Item(434, 1).Interior.Color = 49407
$bands.Cells.Item(520, 1).Interior.Color = 49407

# --- 1d. Rows that had their highlight cleared back to white --------------
# (pre-seed with the already-allocated orange fill so the engine's
# Interior.ThemeColor setter doesn't allocate a throw-away transient fill)
$whiteRows = 144,145,402
foreach ($r in $whiteRows) {
    $wcell = $bands.Cells.Item($r, 1)
    $wcell.Interior.Color = 49407
    $wcell.Interior.ThemeColor = 2
}

# --- 1f. Page setup (A4, portrait) -----------------------------------------
$bands.PageSetup.PaperSize = 9
$bands.PageSetup.Orientation = 1

# --- 1e. Selection on "bands" ------------------------------------------
$bands.Activate()
$bands.Range("A146:C147").Select()

# --- 2. New sheet "Arkusz1", positioned right after "bands" ---------------
$sheet2 = $wb.Worksheets.Add($null, $bands)
$sheet2.Name = "Arkusz1"

$sheet2.Columns.Item(2).ColumnWidth = 36.9
$sheet2.Columns.Item(3).ColumnWidth = 42.7

# Rows pulled from "bands" (row#, B text, C text) - same shared strings.
$sheet2.Range("A1").Value = 11
$sheet2.Range("B1").Value = "Wind speed (gust) [m/s]"
$sheet2.Range("C1").Value = '0[-] SFC="Ground or water surface"'

$sheet2.Range("A2").Value = 146
$sheet2.Range("B2").Value = "u-component of wind [m/s]"
$sheet2.Range("C2").Value = '25000[Pa] ISBL="Isobaric surface"'

$sheet2.Range("A3").Value = 147
$sheet2.Range("B3").Value = "v-component of wind [m/s]"
$sheet2.Range("C3").Value = '25000[Pa] ISBL="Isobaric surface"'

$sheet2.Range("A4").Value = 415
$sheet2.Range("B4").Value = "Temperature [C]"
$sheet2.Range("C4").Value = '2[m] HTGL="Specified height level above ground"'

$sheet2.Range("A5").Value = 417
$sheet2.Range("B5").Value = "Dew point temperature [C]"
$sheet2.Range("C5").Value = '2[m] HTGL="Specified height level above ground"'

$sheet2.Range("A6").Value = 420
$sheet2.Range("B6").Value = "u-component of wind [m/s]"
$sheet2.Range("C6").Value = '10[m] HTGL="Specified height level above ground"'

$sheet2.Range("A7").Value = 421
$sheet2.Range("B7").Value = "v-component of wind [m/s]"
$sheet2.Range("C7").Value = '10[m] HTGL="Specified height level above ground"'

$sheet2.Range("A8").Value = 424
$sheet2.Range("B8").Value = "Precipitation rate [kg/(m^2 s)]"
$sheet2.Range("C8").Value = '0[-] SFC="Ground or water surface"'

$sheet2.Range("A9").Value = 432
$sheet2.Range("B9").Value = "Surface Lifted Index [C]"
$sheet2.Range("C9").Value = '0[-] SFC="Ground or water surface"'

$sheet2.Range("A10").Value = 433
$sheet2.Range("B10").Value = "Convective available potential energy [J/kg]"
$sheet2.Range("C10").Value = '0[-] SFC="Ground or water surface"'

$sheet2.Range("A11").Value = 434
$sheet2.Range("B11").Value = "Convective inhibition [J/kg]"
$sheet2.Range("C11").Value = '0[-] SFC="Ground or water surface"'

$sheet2.Range("A12").Value = 520
$sheet2.Range("B12").Value = "Pressure reduced to MSL [Pa]"
$sheet2.Range("C12").Value = '0[-] MSL="Mean sea level"'

# Same orange highlight for column A, row 10 (A433-equivalent) keeps the
# "applied font" quirk too, the rest are plain orange.
for ($r = 1; $r -le 12; $r++) {
    $cell = $sheet2.Cells.Item($r, 1)
    $cell.Interior.Color = 49407
    if ($r -eq 10) {
        $cell.Font.ThemeColor = 1
    }
}

# --- 2d/2e. Selection + make this the active sheet ------------------------
$sheet2.Activate()
$sheet2.Range("C18").Select()
